# The workbook tracked a "Stem" parameter-count column (F) and a "Size" column (E)
# that were never really meaningful together; the author consolidated on a single
# "Parameters" column (the old Stem numbers), dropped the old "Size" column, and
# swapped the "EfficientNet-B0" benchmark row for a "MobileNetV3(large)" one
# (new model, no accuracy/loss numbers yet). Rows were also re-sorted by the new
# Parameters column, and two blank rows were added at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old "Size" column (E) - old "Stem" column (F) will shift left to
#    become the new column E, carrying its values and per-cell styles with it.
$ws.Range("E1:E11").Delete()

# 2) Re-point the header of what is now column E ("Stem" -> "Parameters"), and
#    let the column re-fit its new (wider) header text.
$ws.Range("E1").Value = "Parameters"
$ws.Columns("E").ColumnWidth = 10.32

# 3) Rewrite rows 2-11 in the new (ascending-by-Parameters) order. Values for
#    A/B/C/D are simply re-typed in place; column E already holds the right
#    number after the shift above, so we only touch it where a row's data moved.
$ws.Range("A2").Value = "MobileNet(alpha=0.25)"
$ws.Range("B2").Value = 0.730934798717498
$ws.Range("C2").Value = 0.72554349899291903
$ws.Range("D2").Value = 224
$ws.Range("E2").Value = "0.2M"

$ws.Range("A3").Value = "MobileNetV2(alpha=0.35)"
$ws.Range("B3").Value = 0.40091994404792702
$ws.Range("C3").Value = 0.86865943670272805
$ws.Range("D3").Value = 224
$ws.Range("E3").Value = "0.4M"

$ws.Range("A4").Value = "MobileNetV2(alpha=0.50)"
$ws.Range("B4").Value = 0.32090279459953303
$ws.Range("C4").Value = 0.89492756128311102
$ws.Range("D4").Value = 224
$ws.Range("E4").Value = "0.7M"

$ws.Range("A5").Value = "MobileNet(alpha=0.50)"
$ws.Range("B5").Value = 0.52614372968673695
$ws.Range("C5").Value = 0.8125
$ws.Range("D5").Value = 224
$ws.Range("E5").Value = "0.8M"

$ws.Range("A6").Value = "MobileNetV3(small)"
$ws.Range("B6").Value = 0.27091637253761203
$ws.Range("C6").Value = 0.90670287609100297
$ws.Range("D6").Value = 224
$ws.Range("E6").Value = "0.9M"

$ws.Range("A7").Value = "MobileNetV2(alpha=0.75)"
$ws.Range("B7").Value = 0.295252114534378
$ws.Range("C7").Value = 0.90851449966430597
$ws.Range("D7").Value = 224
$ws.Range("E7").Value = "1.4M"

$ws.Range("A8").Value = "MobileNet(alpha=0.75)"
$ws.Range("B8").Value = 0.48250153660774198
$ws.Range("C8").Value = 0.816123187541961
$ws.Range("D8").Value = 224
$ws.Range("E8").Value = "1.8M"

$ws.Range("A9").Value = "MobileNetV2(alpha=1.0)"
$ws.Range("B9").Value = 0.29797556996345498
$ws.Range("C9").Value = 0.90760868787765503
$ws.Range("D9").Value = 224
$ws.Range("E9").Value = "2.3M"

# Row 10 used to be "EfficientNet-B0" - the model being benchmarked changed to
# "MobileNetV3(large)", which has no accuracy/loss numbers recorded yet.
$ws.Range("A10").Value = "MobileNetV3(large)"
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 224
$ws.Range("E10").Value = "3.0M"

$ws.Range("A11").Value = "MobileNet(alpha=1.0)"
$ws.Range("B11").Value = 0.41326564550399703
$ws.Range("C11").Value = 0.85778987407684304
$ws.Range("D11").Value = 224
$ws.Range("E11").Value = "3.2M"

# 4) Two trailing blank rows were added below the table (same left/vcenter
#    style as the rest of the sheet, but with no content).
$ws.Range("A12").HorizontalAlignment = -4131
$ws.Range("D12").HorizontalAlignment = -4131
$ws.Range("E12").HorizontalAlignment = -4131

$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("B13").NumberFormat = "0.000_ "
$ws.Range("B13").HorizontalAlignment = -4131
$ws.Range("C13").NumberFormat = "0.000_ "
$ws.Range("C13").HorizontalAlignment = -4131
$ws.Range("D13").HorizontalAlignment = -4131
$ws.Range("E13").HorizontalAlignment = -4131

# 5) Restore the selection the author left the sheet on.
$ws.Range("E11").Select()
